$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text so numeric-looking strings such as
# "213.89" or "1.635.22" are not coerced into numbers by Excel's
# automatic type detection on assignment.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "25.859.10"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.635.22"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -1.07%  "
$ws.Range("D5").Value = "213.89"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("D8").Value = "0.2558"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "0.06357"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "0.07774"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.242"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.636.78"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "1.860.75"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "0.5398"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "0.0₅7856"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "64.30"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "25.865.00"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").Value = "195.89"
$ws.Range("E20").Value = "  -4.34%  "
$ws.Range("D21").Value = "4.354"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "9.876"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "5.948"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("D25").Value = "1.880"
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("D26").Value = "139.31"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "0.1132"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("D28").Value = "6.803"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "15.64"
$ws.Range("D30").Value = "1.232"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("E31").Value = "  -4.08%  "
$ws.Range("D32").Value = "3.239"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "3.163"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "1.525"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "2.350"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "0.8840"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "2.600"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "0.5499"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").Value = "1.121.60"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "0.01552"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "0.9989"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").Value = "5.643"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "0.8095"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "98.95"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").Value = "0.0₈121"
$ws.Range("E45").Value = "  +7.68%  "
$ws.Range("D46").Value = "1.773.24"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "0.4520"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").Value = "54.88"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "0.05039"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  -1.27%  "

# Restore the default (Normal) style on column D so no stray number
# format is left behind now that the text values are committed.
$priceRange.Style = "Normal"
